# Update "想去人数" (F) and "最低票价" (G) values on the "展览" and
# "全部类型" sheets, which carry duplicate data.
#
# Row -> (new F value, new G value or $null if unchanged)
$updates = @{
    2  = @{ F = 2090; G = $null }
    3  = @{ F = 133;  G = $null }
    4  = @{ F = 45;   G = $null }
    6  = @{ F = 1747; G = $null }
    7  = @{ F = 29;   G = $null }
    8  = @{ F = 730;  G = $null }
    9  = @{ F = 368;  G = $null }
    11 = @{ F = 31;   G = $null }
    12 = @{ F = 105;  G = $null }
    15 = @{ F = 145;  G = $null }
    16 = @{ F = 120;  G = $null }
    18 = @{ F = 4060; G = $null }
    19 = @{ F = 11;   G = $null }
    21 = @{ F = 451;  G = $null }
    22 = @{ F = 389;  G = 50 }
    23 = @{ F = 952;  G = $null }
    24 = @{ F = 891;  G = $null }
    26 = @{ F = 24;   G = $null }
    28 = @{ F = 1826; G = $null }
    29 = @{ F = 42;   G = $null }
    30 = @{ F = 37;   G = $null }
    31 = @{ F = 77;   G = $null }
    32 = @{ F = 183;  G = $null }
    33 = @{ F = 18;   G = $null }
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]

        $ws.Range("F$row").Value = $vals.F

        if ($null -ne $vals.G) {
            $ws.Range("G$row").Value = $vals.G
        }
    }
}
